# Iron Mountain Incorporated.xlsx -- "fixed issues restulted from extra files"
#
# The sheet originally mixed in price-history rows that actually belonged to
# several *other* tickers (2337 TT, 005930 KS, CCC LN, UBI FP, IFX GR, XRO AU,
# AMP IM, 6701 JP, 2353 TT, 5803 JP, TIETO FH, 2344 TT, 3665 TT, WKL NA,
# 3661 TT, EXPN LN, 6954 JP, ENR GR, 3443 TT, 6526 JP, 6588 JP) -- an artifact
# of extra files being pulled into the per-ticker workbook. This restores the
# open/close/high/low price columns (D:G) and shares_outstanding (H) to the
# real Iron Mountain (IRM) figures, and re-points the fixed_ticker column (I)
# -- which previously held one of those foreign ticker strings per row -- back
# to "IRM" for every data row. Once nothing references the stray ticker
# strings any more they drop out of the shared-string table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
$data = @(
    @(2,  20.0446039456289,   18.9407787322998,    21.13744765784818,  18.83094486704104,  295348225),
    @(3,  17.2532324070162,   16.75152206420898,   17.97792347044282,  15.7982745393973,    295348225),
    @(4,  17.6141052132338,   17.37033081054688,   18.33975804209759,  16.950812498533,     295348225),
    @(5,  15.45724950684429,  15.88405513763428,   15.94749832592341,  13.63467833245808,   295348225),
    @(6,  19.49847309005985,  21.41549110412598,   21.46825314988352,  18.83015533469323,   295348225),
    @(7,  23.71669827790753,  24.48309516906738,   24.61973920708037,  23.30676616386857,   295348225),
    @(8,  22.5610830668137,   20.28753280639648,   22.5610830668137,   19.57779977266682,   295348225),
    @(9,  20.05731240545477,  21.8851490020752,    22.1785808938808,   19.72720065267744,   295348225),
    @(10, 22.08710231211101,  21.59627723693848,   22.80159463014999,  21.47823168989646,   295348225),
    @(11, 21.87875785647262,  23.00268745422364,   23.18580028335986,  20.78639644607079,   295348225),
    @(12, 24.86245901952425,  25.61159706115723,   26.3607351027902,   24.01727406939481,   295348225),
    @(13, 24.57992771614078,  22.77869987487793,   24.67096409040844,  21.82931699957922,   295348225),
    @(14, 21.75927963688583,  22.46759796142578,   23.01704251657992,  21.28265580347533,   295348225),
    @(15, 23.50203415467352,  23.6434497833252,    24.8825255581116,   22.22928579023469,   295348225),
    @(16, 23.69691008199242,  20.95212173461914,   23.90225540737147,  20.68517124496189,   295348225),
    @(17, 22.29328889142747,  25.91594886779785,   25.93684797571208,  22.00068942157613,   295348225),
    @(18, 25.14634234107422,  23.02659034729004,   25.98289973084963,  22.65793770904856,   295348225),
    @(19, 22.83045914249054,  21.25494766235352,   23.18458903315062,  21.16099576257181,   295348225),
    @(20, 23.71746105764698,  24.13691329956055,   25.34375854346379,  23.27593186483912,   295348225),
    @(21, 23.92713436942469,  23.70961570739746,   24.08464916594844,  22.65202146873482,   295348225),
    @(22, 17.63904590095471,  18.5682258605957,    21.81651311862556,  16.12625054388621,   295348225),
    @(23, 20.64367338234148,  22.1440315246582,    22.34041349614757,  20.09380326286104,   295348225),
    @(24, 21.53049713896848,  20.9047966003418,    23.01452946735444,  20.66414337088933,   295348225),
    @(25, 24.56592209633205,  27.57115173339844,   33.83546311745052,  22.69891145474148,   295348225),
    @(26, 30.99653985381529,  33.41163635253906,   34.07786975988178,  30.46355312794111,   295348225),
    @(27, 35.87349989160018,  36.93704223632812,   38.02590924111247,  35.43457668300683,   295348225),
    @(28, 37.33379984092087,  39.03584289550781,   39.97667368221356,  36.09361573186261,   295348225),
    @(29, 45.44514680828665,  39.77207946777344,   45.47112925344577,  36.42020857934622,   295348225),
    @(30, 48.73578779424142,  47.11359786987305,   51.39266747505702,  46.99083823834282,   295348225),
    @(31, 43.2299190760566,   43.07004165649414,   44.24249865535388,  39.55266727160468,   295348225),
    @(32, 40.15734962697366,  44.98162078857422,   45.37690442138325,  39.42966488893347,   295348225),
    @(33, 45.67030079147715,  49.5859375,          50.08561186568342,  44.44382356772381,   295348225),
    @(34, 48.50248646405817,  50.78236389160156,   51.04895943567091,  46.9672501835215,    295348225),
    @(35, 52.48008746235487,  57.07186508178711,   58.44753824296802,  52.31277552625755,   295348225),
    @(36, 55.69396814178256,  55.47795486450195,   58.53971637298933,  53.07362709449612,   295348225),
    @(37, 66.28448134583748,  64.02758026123047,   66.43620171514202,  61.33448040183399,   295348225),
    @(38, 76.6834734775279,   74.11174011230469,   76.845997479876,    70.29716713070502,   295348225),
    @(39, 86.00899313504182,  98.78031158447266,   99.59899206178066,  85.30589289350355,   295348225),
    @(40, 115.3187905390577,  119.902473449707,    126.2110917081674,  112.5666481624115,   295348225),
    @(41, 102.6695405329658,  99.06094360351562,   109.4088482530502,  94.7891442315747,    295348225),
    @(42, 86.90635032474394,  88.25472259521484,   88.72714998426342,  71.18840612535941,   295348225),
    @(43, 101.1917197881852,  96.56955718994141,   102.1637667235841,  95.54792073349807,   295348225)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 4).Value  = $entry[1]   # D: open_price
    $ws.Cells.Item($r, 5).Value  = $entry[2]   # E: close_price
    $ws.Cells.Item($r, 6).Value  = $entry[3]   # F: high_price
    $ws.Cells.Item($r, 7).Value  = $entry[4]   # G: low_price
    $ws.Cells.Item($r, 8).Value  = $entry[5]   # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value  = "IRM"       # I: fixed_ticker -> back to IRM for every row
}
